# Update "想去人数" (want-to-go count) figures on both the "展览" and
# "全部类型" worksheets, which hold identical data in this workbook.

$wb = $excel.ActiveWorkbook

$updates = @{
    "F2"  = 163
    "F3"  = 1750
    "F4"  = 801
    "F6"  = 1129
    "F8"  = 12116
    "F15" = 13534
    "F16" = 13591
    "F21" = 1001
    "F24" = 2053
    "F25" = 191
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($addr in $updates.Keys) {
        $ws.Range($addr).Value = $updates[$addr]
    }
}
